# Fixed setting end dates for viranomaispaatos and kompostori
#
# Adds a new decision row (row 3) to Sheet1 that mirrors the existing
# decision row (row 2), but with an updated decision date, recipient and
# "voimassa alkaen" (valid-from) date - matching the additional test-data
# row that was introduced together with the fix for end-date handling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2 (values + formatting) into the new row 3.
$ws.Range("A2:P2").Copy($ws.Range("A3:P3")) | Out-Null

# Column C (Paatospaivamaara / decision date) - update to new date, keep as text.
$ws.Range("C3").Value = "20.6.2022"

# Column G (Lahettaja/vastaanottaja) - new recipient.
$ws.Range("G3").Value = "vastaanottaja Karita Pyykoski"

# Column K (Voimassa alkaen 1 / valid-from date) - update to new date.
# Force a text number format first so the day/month value (1.7.2022) is not
# auto-converted into a date serial number, matching how the original data
# stores these as plain text values.
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = "1.7.2022"
$ws.Range("K3").NumberFormat = "General"

# Match the selection left behind in the authored workbook.
$ws.Range("G3").Select() | Out-Null
